$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.13"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "21"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.16"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "21"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.370"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "21"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "21"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.387"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "21"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.383"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "21"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8131"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "21"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9611"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "21"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1424"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "21"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03552"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07365"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "21"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03035"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "21"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.458"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "21"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09391"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "21"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001590"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "21"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04839"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "21"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "One"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005893"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "21"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006281"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "21"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004087"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "21"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009867"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "21"

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "21"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.686"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "21"

$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "21"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3253"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "21"

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "21"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "21"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "21"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "21"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "21"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "21"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "21"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "21"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "21"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "21"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "21"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "21"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "21"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "21"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03855"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "21"

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1076"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "21"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002441"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "21"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003031"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "21"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005749"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "21"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005650"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "21"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "21"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "21"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07683"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "21"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "21"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "21"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "21"
